$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers (unchanged values, kept for completeness) ---
$ws.Range("A1").Value = "Tiles"
$ws.Range("F1").Value = "NPCs"
$ws.Range("K1").Value = "Key"
$ws.Range("N1").Value = "Map Construction"

# --- Row 3 headers ---
$ws.Range("A3").Value = "Char"
$ws.Range("B3").Value = "Name"
$ws.Range("C3").Value = "Colour"
$ws.Range("D3").Value = "Type"
$ws.Range("F3").Value = "Char"
$ws.Range("G3").Value = "Name"
$ws.Range("H3").Value = "Colour"
$ws.Range("I3").Value = "Type"
$ws.Range("K3").Value = "Type #"
$ws.Range("L3").Value = "Represents"
$ws.Range("N3").Value = "Format"
$ws.Range("O3").Value = "Min"
$ws.Range("P3").Value = "Max"

# --- Tiles table (A:D) ---
$ws.Range("A4").Value = "H"
$ws.Range("B4").Value = "Player"
$ws.Range("C4").Value = "Colour.WHITE"
$ws.Range("D4").Value = 1

$ws.Range("A5").Value = [char]0x2592
$ws.Range("B5").Value = "Ascii 177"
$ws.Range("C5").Value = "Colour.GREY_40"
$ws.Range("D5").Value = 1

$ws.Range("A6").Value = "n"
$ws.Range("B6").Value = "Grass"
$ws.Range("C6").Value = "Colour.GREEN_BB"
$ws.Range("D6").Value = 0

$ws.Range("A7").Value = "w"
$ws.Range("B7").Value = "Water"
$ws.Range("C7").Value = "Colour.BLUE_BB"
$ws.Range("D7").Value = 2

$ws.Range("A8").Value = "E"
$ws.Range("B8").Value = "Bank"
$ws.Range("C8").Value = "Colour.SADDLE_BROWN"
$ws.Range("D8").Value = 2

$ws.Range("A9").Value = "B"
$ws.Range("B9").Value = "Brick/Wall"
$ws.Range("C9").Value = "Colour.RED_BB"
$ws.Range("D9").Value = 1

$ws.Range("A10").Value = "c"
$ws.Range("B10").Value = "Pathway"
$ws.Range("C10").Value = "Colour.GREY_70"
$ws.Range("D10").Value = 0

$ws.Range("A11").Value = "@"
$ws.Range("B11").Value = "Bridge"
$ws.Range("C11").Value = "Colour.GREY_D0"
$ws.Range("D11").Value = 0

# --- NPCs table (F:I) ---
$ws.Range("F4").Value = "S"
$ws.Range("G4").Value = "Sheep"
$ws.Range("H4").Value = "Colour.WHITE"
$ws.Range("I4").Value = 3

$ws.Range("F5").Value = "D"
$ws.Range("G5").Value = "Duck"
$ws.Range("H5").Value = "Colour.YELLOW"
$ws.Range("I5").Value = 3

$ws.Range("F6").Value = "M"
$ws.Range("G6").Value = "Cow"
$ws.Range("H6").Value = "Colour.BLACK"
$ws.Range("I6").Value = 3

# --- Key table (K:L) ---
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = "Anything can pass."
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = "Nothing can pass."
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = "Sheep/Cow can't pass."
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = "Is an NPC."

# --- Map Construction table (N:P) ---
$ws.Range("N4").Value = "Name (alpha-numeric)"
$ws.Range("O4").Value = "n/a"
$ws.Range("P4").Value = "n/a"

$ws.Range("N5").Value = "Width"
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 60

$ws.Range("N6").Value = "Height"
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 60

$ws.Range("N7").Value = "Start X Position"
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 59

$ws.Range("N8").Value = "Start Y Position"
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 59

$ws.Range("N9").Value = "NPC Count"
$ws.Range("O9").Value = 0
$ws.Range("P9").Formula = "=(60*60)-1"

$ws.Range("N10").Value = "Map data"
$ws.Range("O10").Value = "n/a"
$ws.Range("P10").Value = "n/a"

# --- Column L width (Format/Represents column got wider to fit new text) ---
$ws.Columns("L").ColumnWidth = 20.166666666666668

# --- View: clear the scrolled-right view and move selection to K9 ---
$ws.Range("K9").Select()
